$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-arrange header row (Category, Product, L, M, S, XL, XS)
$ws.Range("C1").Value = "L"
$ws.Range("D1").Value = "M"
$ws.Range("E1").Value = "S"
$ws.Range("F1").Value = "XL"
$ws.Range("G1").Value = "XS"

# Row 2 data (All / Prod1), written per new query result
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 35
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 30
$ws.Range("G2").Value = 22

# Remove the stray extra columns (H:P) that were left over on row 2
$ws.Range("H2:P2").ClearContents()

# Row 3 data (All / Prod2), written per new query result
$ws.Range("C3").Value = 25
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 24

# Row 3 no longer has XL / XS figures
$ws.Range("F3:G3").ClearContents()
